$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.867.07"
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("D3").Value = "2.263.66"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'300.92"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'100.77"
$ws.Range("E6").Value = "  +6.54%  "
$ws.Range("D7").Value = "'0.558"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.505"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "'35.42"
$ws.Range("E10").Value = "  +3.77%  "
$ws.Range("D11").Value = "'0.0774"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "2.606.55"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "2.258.95"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'13.53"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "46.833.01"
$ws.Range("E17").Value = "  +4.43%  "
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "'12.74"
$ws.Range("E19").Value = "  -5.84%  "
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").Value = "'65.35"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'247.18"
$ws.Range("E23").Value = "  +3.55%  "
$ws.Range("D24").Value = "'2.80"
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "'1.86"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").Value = "'41.99"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").Value = "'20.10"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "'2.83"
$ws.Range("E31").Value = "  +10.59%  "
$ws.Range("D32").Value = "'145.37"
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("D33").Value = "'5.37"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").Value = "'3.24"
$ws.Range("E34").Value = "  +11.80%  "
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  +10.92%  "
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").Value = "'16.07"
$ws.Range("E38").Value = "  +18.41%  "
$ws.Range("E39").Value = "  -4.25%  "
$ws.Range("D40").Value = "'3.86"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("D42").Value = "'3.11"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'1.91"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "'91.56"
$ws.Range("E45").Value = "  +20.01%  "
$ws.Range("D46").Value = "1.778.06"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "'71.26"
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("D50").Value = "'7.80"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'93.53"
$ws.Range("E51").Value = "  -1.79%  "
